$wb = $excel.ActiveWorkbook

# --- Week 17: append new per-play / per-game logged data to the long run-length shared strings ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "5 3 1 7 2 2 2 1 -3 -5 13 7 4 10 3 2 3 2 3 1 3 9 6 4 6 1 2 1 5 6 5 6 9 5 4 9 3 3 5 2 4 2 15 0 11 14 0 1 4 9 12 3 2 3 3 4 24 1 1 19 5 2 2 8 4 -6 3 -3 4 7 0 6 9 7 7 2 2 2 1 2 2 0 6 0 1 8 5 8 11 3 -1 8 9 0 14 4 8 7 2 2 5 -3 4 12 0 1 0 -2 2 1 74 1 3 3 1 -2 4 6 4 -3 0 0 2 0 2 4 2 5 1 20 3 1 26 5 74 5 -1 40 5 1 12 5 8 1 0 1 1 2 3 7 5 3 10 3 2 0 2 1 2 0 4 6 12 4 19 9 7 6 17 -4 1 7 4 11 3 2 4 2 2 4 0 7 -2 9 0 4 4 7 0 5 9 11 2 7 2 4 5 -1 5 3 14 56 4 0 14 5 13 -1 4 1 -2 0 0 9 11 5 3 4 5 11 1 2 7 3 6 6 11 9 3 -1 1 -1 7 1 0 2 1 0 1 8 0 20 5 13 6 6 2 1 1 0 7 6 3 3 8 8 -2 2 11 11 2 3 3 3 12 3 -4 5 1 3 3 10 13 2 5 5 3 3 5 2 4 9 7 1 1 1 15 82 24 16 -1 6 2 7 5 3 5 3 2 9 19 7 5 1 1 3 2 8 1 -3 1 2 10 12 11 10 4 3 0 4 5 -2 3 4 2 2 9 4 7 5 1 4 8 3 7 5 11 12 2 6 1 4 7 -1 2 2 7 20 1 -4 8 6 12 6 1 3 5 2 18 3 5 3 2 6 2 0 0 0 9 1 8 1 6 0 3 6 7 4 3 20 1 4 3 6 3 1 1 5 0 4 13 7 4 6 6 3 4 11 4 9 23 1 4 14 18 7 3 3 1 6 2 8 -6 -1 5 -3 9 -1 6 6 2 2 0 5 10 5 2 6 6 12 -2 9 -2 -3 1 7 9 7 6 6 13 3 3 0 8 9 27 1 -2 24 2 3 3 1 3 7 9 8 2 1 2 5 11 12 5 12 6 3 3 -4 7 2 2 7 4 2 13 10 6 4 1 3 0 1 4 4 1 4 14 3 6 1 6 7 18 -3 12 1 -2 11 1 1 -1 5 6 2 5 -2 9 11 23 14 7 3 -1 2 3 5 9 7 1 5 7 5 2 7 1 4 3 4 22 5 8 5 0 6 0 0 6 4 1 16 4 2 1 5 5 7 11 3 10 5 1 1 8 6 10 4 3 19 1 12 3 14 4 -2 20 4 7 3 10 4 3 21 4 3 5 2 2 3 0 3 8 18 4 2 5 4 2 1 4 6 8 5 4 5 3 5 6 2 3 10 2 2 3 3 0 5 2 12 2 1 10 2 4 6 0 11 10 17 -2 0 1 4 3 -8 7 4 2 5 11 31 0 3 4 3 0 5 -1 5 21 -2 15 5 25 1 7 5 3 1 2 6 1 23 1 2 3 4 -3 4 1 25 4 4 20 4 3 3 2 4 4 18 9 4 13 4 4 4 5 6 0 2 5 6 0 4 -3 7 5 10 6 2 -6 3 -1 7 2 1 14 1 3 24 -1 1 7 4 5 3 2 11 5 7 2 14 18 6 1 27 3 5 12 5 4 13 5 13 12 1 1 6 3 4 3 4 0 5 9 2 2 1 4 9 7 18 9 4 1 3 0 7 6 2 2 6 0 0 1 1 3 -1 3 5 -1 8 34 1 4 2 2 7 12 3 2 0 4 5 4 4 4 6 5 8 1 4 19 0 9 4 4 10 3 1 8 8 38 1 9 12 7 -3 4 3 15 4 2 3 4 4 0 8 4 -3 1 11 0 3 15 5 0 13 4 1 1 2 5 5 2 3 9 1 4 4 6 6 6 4 5 3 3 3 6 0 0 6 6 -3 7 4 3 0 3 -3 22 3 6 4 2 0 1 14 13 5 4 2 8 0 1 7 4 2 3"
$ws.Range("C2").Value = "2 -3 8 0 2 5 7 1 20 -3 4 -2 0 1 0 1 0 3 -1 1 5 3 1 3 1 0 3 1 2 1 0 4 4 5 5 3 -1 0 5 15 3 3 19 9 5 1 -1 5 -1 3 -1 4 9 3 40 2 2 11 5 5 10 1 2 6 3 4 11 2 5 3 2 1 2 4 0 1 0 9 0 -1 6 9 5 15 10 8 -5 38 3 4 -1 5 3 3 0 0 9 4 11 2 3 1 0 1 4 1 1 3 1 3 4 -1 2 0 25 2 2 0 0 58 7 0 1 1 0 -1 5 -1 8 1 1 1 8 2 0 1 -2 9 4 -2 2 3 7 2 7 5 -2 0 10 16 5 6 2 1 -2 24 -1 4 6 1 37 0 5 -4 3 1 10 5 21 7 7 10 3 2 3 2 -1 80 1 -2 2 4 0 19 -4 1 14 10 3 0 7 19 3 5 11 4 -3 -3 4 1 -1 -3 3 9 6 8 7 2 9 2 10 0 11 0 0 4 2 14 6 6 -4 5 9 6 17 2 34 2 6 2 4 10 5 4 0 0 1 6 3 3 4 0 8 1 1 10 1 9 6 1 0 4 -3 6 2 -4 1 0 0 9 3 -1 2 -1 3 -2 3 2 6 -1 3 3 2 -9 9 -3 7 11 1 4 54 2 5 3 0 10 6 2 0 4 2 -1 0 4 1 6 3 1 2 1 1 1 -3 0 7 8 16 2 9 4 0 5 10 -1 -2 -2 -3 -2 1 4 4 1 5 3 6 1 1 1 3 6 7 3 18 3 3 8 -4 1 3 77 1 5 4 2 3 2 2 3 16 3 7 0 2 9 2 4 5 1 3 7 11 10 2 1 6 5 -1 1 8 7 -1 3 3 2 5 -4 5 10 -1 9 8 2 2 16 3 0 1 -1 2 5 0 4 1 5 6 3 1 -1 4 1 9 5 3 -1 0 3 5 5 2 5 1 7 1 3 14 9 3 -1 2 4 31 0 4 -1 19 3 3 1 6 18 0 -6 -3 11 0 4 3 0 3 1 3 12 2 5 2 13 0 4 9 5 1 3 12 4 3 -2 3 10 11 14 1 8 7 12 3 2 11 3 1 4 3 -1 0 0 -1 3 8 5 -4 -1 3 3 0 1 2 5 2 7 4 5 2 2 10 2 3 2 7 4 21 2 0 1 3 9 3 0 8 -1 8 1 1 3 5 9 1 4 4 8 11 7 5 -1 6 5 7 6 4 9 5 8 1 2 10 3 0 0 13 2 -1 5 4 8 3 1 5 4 1 5 -6 5 4 5 4 5 3 12 2 15 24 17 1 1 2 9 3 10 10 1 8 1 2 4 1 12 9 16 3 3 5 4 5 1 5 6 0 26 3 1 5 0 3 3 -4 1 2 9 -5 8 2 3 11 1 8 6 11 3 5 -1 2 5 4 6 2 4 0 9 4 3 6 0 2 3 2 2 20 5 6 2 7 0 1 9 2 2 0 1 -1 3 1 5 1 2 3 12 9 4 3 9 4 8 -3 6 11 11 4 1 5 1 2 7 7 5 1 3 1 1 4 1 1 1 -2 5 4 2 2 -2 5 0 2 7 8 3 1 3 7 3 0 2 2 0 9 -3 0 1 1 5 -1 1 14 8 10 2 1 0 16 1 3 9 11 2 15 7 1 1 4 20 -4 8 9 1 2 6 -1 4 2 14 9 1 -2 4 3 1 9 1 -1 28 3 1 17 7 3 1 5 0 2 14 2 2 1 1 32 6 -1 -3 7 -4 3 4 13 -1 1 4 3 5 0 1 -1 2 1 3 17 9 1 0 1 15 2 2 3 2 3 8 17 0 2 6 5 2 4 1 3 3 5 4 -1 -1 8 -6 2 30 0 1 0 13 1 5 5 3 4 2 1 6 6 -2 -1 -4 10 3 7 1 6 5 0 5 0 2 2 4 4 2 7 3 3 11 3 6 1 1 3 4 4 5 8 1 4 3 8 0 4 10 1 8 3"
$ws.Range("B3").Value = "16 19 12 5 7 55 -2 0 11 7 10 34 -4 12 23 6 3 27 6 -6 10 4 7 8 12 12 14 4 2 8 16 5 8 12 10 8 5 5 6 9 12 9 25 3 4 8 4 22 8 11 11 7 9 7 12 2 6 3 0 12 16 3 10 8 29 9 10 0 3 11 3 -2 11 4 5 0 4 30 6 -4 6 20 28 7 8 4 5 15 13 23 9 3 1 2 2 17 17 17 42 5 13 7 9 18 6 20 9 37 31 11 8 2 18 19 18 13 11 4 9 8 -6 9 8 3 19 13 8 6 4 50 6 11 3 5 9 10 10 10 1 16 10 8 9 18 15 13 12 8 16 14 6 10 1 15 16 7 3 6 11 40 59 3 11 30 18 8 32 6 2 18 6 3 6 4 6 4 9 2 15 9 2 8 8 5 22 19 15 6 1 6 5 16 8 16 5 9 13 7 6 12 22 16 32 3 1 4 12 19 8 29 4 3 9 5 7 7 14 7 4 3 1 6 10 8 2 3 3 9 9 17 16 0 5 3 7 3 5 11 10 6 27 -1 8 17 33 2 7 11 8 10 41 34 6 20 17 9 2 7 5 3 39 19 0 6 13 11 6 6 6 6 14 5 26 -6 10 24 16 32 7 12 42 5 27 10 4 23 3 12 -6 39 5 14 18 3 5 81 19 5 17 23 10 15 43 10 8 9 22 17 6 10 14 9 8 15 9 7 4 15 30 7 0 5 10 5 6 6 10 11 2 18 13 28 8 4 19 3 9 -3 1 9 6 10 12 9 -1 5 8 25 7 19 9 0 23 4 11 4 9 1 5 26 91 6 10 9 14 38 5 5 9 19 7 19 28 41 27 0 3 4 1 5 12 5 7 24 19 16 12 15 -4 9 13 22 21 12 10 6 3 5 15 3 7 14 18 6 37 12 1 21 8 7 6 4 14 14 9 18 19 17 9 9 15 -7 6 -5 3 10 8 8 2 7 2 0 1 15 17 7 11 8 13 5 53 14 20 23 5 6 4 16 2 1 6 -3 24 25 6 24 13 8 -3 20 27 25 15 20 7 5 6 17 10 8 7 6 7 16 0 15 13 17 19 43 15 2 18 27 7 23 19 21 -6 19 4 28 16 18 3 4 6 24 36 8 6 4 5 9 11 12 14 2 8 4 16 14 33 -2 4 14 6 14 14 -1 23 14 0 19 2 0 13 4 21 6 10 8 7 14 11 16 9 36 12 10 12 5 25 13 16 12 22 -2 2 28 9 7 7 2 1 5 -1 19 28 16 45 6 5 8 17 21 20 5 9 7 9 7 34 19 17 6 2 3 23 13 46 1 9 1 6 39 4 3 16 10 5 30 5 12 14 5 19 18 5 8 8 18 13 -3 27 12 13 10"
$ws.Range("C3").Value = "16 2 14 21 11 7 6 6 1 17 19 4 3 8 14 10 19 24 4 10 8 5 19 21 13 12 6 20 3 13 6 21 12 1 13 28 28 8 2 5 3 9 9 15 19 16 5 1 2 7 5 7 13 17 25 4 12 10 1 14 18 42 11 1 8 12 0 11 9 17 7 13 0 17 38 13 5 4 10 12 13 2 5 20 5 10 8 7 5 11 13 4 14 14 2 20 13 5 10 12 2 4 -2 5 7 12 8 4 -1 15 32 11 3 17 6 5 7 5 4 8 11 14 2 11 12 2 2 11 17 7 9 32 9 11 5 18 2 -1 20 10 17 10 9 39 5 7 -4 17 29 9 12 0 5 0 7 11 7 7 12 5 9 2 8 12 13 6 10 14 1 11 2 14 9 8 5 7 6 8 6 6 10 7 10 22 4 11 9 12 3 8 4 17 11 1 4 6 27 38 7 16 7 0 40 -3 13 13 43 21 7 10 5 42 9 18 13 10 12 13 10 8 7 11 18 52 1 6 2 17 3 0 9 19 8 -2 3 31 7 9 24 10 12 9 8 6 2 36 14 5 25 42 1 31 9 13 5 11 7 8 3 -1 5 -6 2 18 5 6 10 2 9 8 2 9 13 11 8 11 1 22 37 12 20 9 1 10 9 2 9 31 20 15 16 6 10 4 42 6 45 2 15 21 6 23 14 9 1 32 16 11 15 8 8 44 12 2 6 20 35 14 7 9 7 21 4 55 6 5 8 16 4 7 19 7 11 52 69 1 14 6 6 1 10 18 16 5 5 6 5 4 5 4 15 6 6 6 13 5 11 6 3 6 16 1 11 12 7 3 7 9 9 5 18 5 12 6 9 7 6 8 7 5 10 1 8 7 6 11 17 5 8 5 6 2 40 11 9 4 6 6 6 0 15 8 8 8 44 7 8 4 12 19 9 8 6 9 6 15 15 13 11 7 5 22 11 9 36 1 9 2 19 32 4 6 4 9 2 15 7 6 3 14 7 10 5 12 11 44 6 13 14 5 6 12 11 -1 8 5 7 0 7 5 24 2 8 9 4 7 25 9 22 16 10 2 8 8 11 3 9 23 -4 10 12 14 18 10 3 6 -4 5 4 4 9 15 0 8 5 4 8 10 2 27 10 12 4 43 -1 16 2 8 13 0 6 18 18 8 10 10 6 20 7 6 5 4 21 7 9 10 24 5 -1 29 3 18 7 5 8 -1 5 7 7 2 10 5 9 18 4 34 0 8 6 11 11 7 9 17 5 19 12 9 17 23 9 23 2 12 15 3 1 4 11 9 16 15 19 49 8 8 -2 7 3 8 18 3 14 4 5 11 -2 7 5 64 4 13 15 26 10 7 4 10 6 15 7 5 -9 0 13 9 14 3 5 7 6 18 19 5 10 9 26 21 -6 17 2 5 2 10 1 20 6 5 22 4 1 20 7 18 19 6 14 9 0 10 17 17 1 18 18 5 4 10 8 8 7 3 10 7 5 10 29 4 19 5 18 14 19 11 -2 15 4 14 7 11 19 7 6 46 5 1 7 5 9 7 5 12 29 5 6 1 13 10 3 -2 -4 8 8 8 -2 7 6 8 2 9 11 5 2 0 4 5 9 6 4 8 2 16 10 24 8 25 12 2 12 12 11 5 6 18 9 4 11 1 -4 2 8 10 4 5 9 12 5 10"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = "61 39 66 71 65 65 66 62 60 58 59 56 67 60 67 58 57 64 62 65 65 65 61 65 60 65 73 23 62 70 69 67 68 66 58 57 57 65 73 73 60 63 65 65 62 66 65 65 61 67 64 62 65 61 66 67"
$ws.Range("D3").Value = "53 51 59 62 43 62 59 45 53 47 56 43 40 61 48 52 50 55 41 51 49 45 43 49 59 49 66 43 49 57 48 54 47 43 40 36 48 46 56 50 33 47 52 42 35 40 43 37 63 57 45 29 53 52 35 42 42 44 54 40 33 42 35 42 47 17 51 27 41 41 49 47 41 41 48 46 54 50 54 32 68 49 41 40 38 39 51 41 64 49 60 16 46 47 39 39 40 55 46 45 52 47 44 44 42 54 53 47 44 53 41 43 42 46 42 38 41 34 39 41 39 51 39 35"
$ws.Range("B5").Value = "19 14 25 24 24 37 28 23 15 15 20 14 17 18 27 30 13 21 22 33 23 22 28 27 20 22 25 4 11 28 22 32 28 39 13 31 26 30 41 29 19 29 39 19 22 20 22 79 28 32 47 28 17 21 34 23"
$ws.Range("D4").Value = "2 0 0 11 4 0 9 2 19 6 10 0 0 0 0 0 0 6 0 6 6 10 0 3 5 10 11 9 14 5 0 0 0 0 4 1 0 10 16 20 0 0 12 0 3 7 2 0 6 12 7 0 12 0 0 0 2 4 0 0 0 0 5 -5 6 0 0 0 0 0 20 0 0 0 7 0 7 0 13 0 0 5 8 0 0 0 0 0 0 22 8 0 9 12 0 4 0 5 0 0 15 11 0 10 0 8 19 0 8 0 14 7 9 0 0 0 0 0 0 3 13 0 0 0"
$ws.Range("B6").Value = "20 25 45 0 15 17 17 14 22 21 16 22 20 31 22 25 20 17 18 40 19 24 27 46 46 16 22 20 15 21 19 15 24 18 22 20 14 44 27 19 26 23 23 24 17 19 18 14 20 13 15 26 21 5 13 31 23 17 20 10 21 0"
$ws.Range("D5").Value = "0 0 8 0 0 0 11 0 0 0 0 0 0 0 1 1 0 8 0 0 0 1 11 5 0 0 0 0 4 0 0 11 0 0 2 9 22 0 0 0 0 0 6 0 10 0 0 7 12 0 0 0 73 0 5 0 0 4 0 0 0 0 0 0 15 15 3 3 0 8 0 0 2 11 0 7 -1 11 2 0 0 0 -4 -3 0 5 0 0 0 22 0 0 0 0 8 6 0 0 0 12 0 10 0 9 12 10 0 8 0 13 4 0 2 3 9 0 0 11 0 1 10 -4 0 0 14 0 39 0 1 0 7"

# --- Week 17: updated cumulative season totals on the summary sheets ---

# OFF sheet (Home = row 2, Road = row 3)
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 389
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 23
$ws.Range("F2").Value = 133
$ws.Range("G2").Value = 123
$ws.Range("I2").Value = 19
$ws.Range("J2").Value = 78
$ws.Range("N2").Value = 52
$ws.Range("O2").Value = 57
$ws.Range("P2").Value = 24
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 341
$ws.Range("F3").Value = 186
$ws.Range("G3").Value = 63
$ws.Range("H3").Value = 51
$ws.Range("I3").Value = 101
$ws.Range("J3").Value = 107
$ws.Range("L3").Value = 551
$ws.Range("M3").Value = 332
$ws.Range("Q3").Value = 1071

# DEF sheet (Home = row 2, Road = row 3)
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 407
$ws.Range("D2").Value = 31
$ws.Range("F2").Value = 116
$ws.Range("G2").Value = 127
$ws.Range("J2").Value = 55
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 43
$ws.Range("B3").Value = 18
$ws.Range("C3").Value = 355
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 200
$ws.Range("H3").Value = 68
$ws.Range("I3").Value = 110
$ws.Range("J3").Value = 106
$ws.Range("L3").Value = 561
$ws.Range("M3").Value = 394
$ws.Range("Q3").Value = 1067

# ST sheet
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 159
$ws.Range("D2").Value = 124
$ws.Range("F2").Value = 182
$ws.Range("G2").Value = 173
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 38
$ws.Range("B3").Value = 103

# TURNS sheet (Road = row 3)
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 18
$ws.Range("E3").Value = 14

# PEN sheet
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 41
$ws.Range("D2").Value = 38
$ws.Range("B3").Value = 25
$ws.Range("D4").Value = 20
